# Adds new nomenclature rows to the "Code / Libellé niveau 1 / Libellé niveau 2 /
# Libellé niveau 3 / Description / Commentaire" table of the document.
#
# Each new row follows the pattern observed for every inserted row in the diff:
#   Col1 (Code)              -> the short code
#   Col2 (Libellé niveau 1)  -> left empty
#   Col3 (Libellé niveau 2)  -> the label text
#   Col4 (Libellé niveau 3)  -> left empty
#   Col5 (Description)       -> left empty
#   Col6 (Commentaire)       -> left empty

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Clean-CellText($text) {
    # Word cell ranges end with a cell-mark (CR + BEL, chars 13/7); strip
    # only that (not general whitespace, since some codes such as "TSU "
    # legitimately keep a trailing space) so comparisons work reliably.
    return $text.TrimEnd([char]13, [char]7)
}

function Find-RowByCode($table, $code) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        if ((Clean-CellText $table.Cell($i, 1).Range.Text) -eq $code) {
            return $i
        }
    }
    return -1
}

function Insert-RowBefore($table, $anchorCode, $newCode, $label) {
    $idx = Find-RowByCode $table $anchorCode
    $table.Rows.Add($table.Rows.Item($idx))
    $table.Cell($idx, 1).Range.Text = $newCode
    $table.Cell($idx, 3).Range.Text = $label
}

function Append-Row($table, $newCode, $label) {
    $table.Rows.Add()
    $idx = $table.Rows.Count
    $table.Cell($idx, 1).Range.Text = $newCode
    $table.Cell($idx, 3).Range.Text = $label
}

# A) Before HOSPIT: MED, PARAMED
Insert-RowBefore $t "HOSPIT" "MED" "Medicalisé"
Insert-RowBefore $t "HOSPIT" "PARAMED" "Paramédicalisé"

# B) Before TSU (code cell literally contains "TSU " with a trailing space):
#    MEDC, PHARMA, INF, MEDSPE, DENT, AUTREPRO
Insert-RowBefore $t "TSU " "MEDC" "Médecin Généraliste"
Insert-RowBefore $t "TSU " "PHARMA" "Pharmacien"
Insert-RowBefore $t "TSU " "INF" "Infirmier"
Insert-RowBefore $t "TSU " "MEDSPE" "Médecin autre spécialité"
Insert-RowBefore $t "TSU " "DENT" "Dentiste"
Insert-RowBefore $t "TSU " "AUTREPRO" "Autre professionnel de santé"

# C) Before AASC: MSP, ISP, SP
Insert-RowBefore $t "AASC" "MSP" "Medecin Sapeur-Pompier"
Insert-RowBefore $t "AASC" "ISP" "Infirmier Sapeur-Pompier"
Insert-RowBefore $t "AASC" "SP" "Secouriste"

# D) Before AUTRE: HELIFSI, VLFSI, FFSI, DGDD
Insert-RowBefore $t "AUTRE" "HELIFSI" "Police Nationale"
Insert-RowBefore $t "AUTRE" "VLFSI" "Gendarmerie Nationale"
Insert-RowBefore $t "AUTRE" "FFSI" "Police Municipale"
Insert-RowBefore $t "AUTRE" "DGDD" "Douane"

# E) Appended at the very end: ADM, DAE, INCONNU
Append-Row $t "ADM" "Institutions administratives et sociales"
Append-Row $t "DAE" "Défibrillateur Automatique"
Append-Row $t "INCONNU" "Autre moyen"
